# Watershed_Stats.xlsx revision: reorder the "% of area" breakdown columns so
# the Bare-Land percentage (old column U) leads the block, then add two new
# summary columns: "% Disturbed" and "% Undisturbed".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fagaalu_Revised")

# --- 1. Reorder columns: cut the old "% Bare Land" column (U) and insert it
#        before column O, pushing O..T right into P..U. Column V (% Water)
#        is untouched by this move.
$ws.Columns("U:U").Cut()
$ws.Columns("O:O").Insert()

# --- 2. Add the two new summary columns after V: "% Disturbed" (W) and
#        "% Undisturbed" (X), plus a trailing spacer column Y.
$ws.Range("W1").Value = "% Disturbed"
$ws.Range("X1").Value = "% Undisturbed"

$ws.Range("W2").Formula = '=SUM(L2,F2:G2)/N2'
$ws.Range("W3:W5").Formula = '=SUM(L3,F3:G3)/N3'
$ws.Range("W6").Formula = '=SUM(L6,F6:H6)/N6'

$ws.Range("X2").Formula = '=SUM(I2:K2)/N2'
$ws.Range("X3:X6").Formula = '=SUM(I3:K3)/N3'
$ws.Range("X6").Formula = '=SUM(I6:K6,M6)/N6'

# Trailing empty helper column.
$ws.Range("Y2:Y6").Value = ""

# --- 3. Formatting to match the new columns.
$ws.Range("W1").Style = $ws.Range("V1").Style
$ws.Range("X1").Style = $ws.Range("V1").Style

$ws.Range("W2:W6").NumberFormat = "0.0%"
$ws.Range("X2:X6").NumberFormat = "0.0%"

$ws.Columns("W:W").ColumnWidth = 14.43
$ws.Columns("X:X").ColumnWidth = 13

# --- 4. View state: refreeze panes further right and move the active cell.
$ws.Range("K2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("O4").Select()
